$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 61
$ws.Range("H61").Value = 440.7143
$ws.Range("I61").Value = 430.83334
$ws.Range("J61").Value = 500
$ws.Range("K61").Value = 1292.50002
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -1120.50002
$ws.Range("N61").Value = -1844

# Row 62
$ws.Range("H62").Value = 4091.7273
$ws.Range("J62").Value = 3882.2
$ws.Range("L62").Value = 3882.2
$ws.Range("N62").Value = -5130.2

# Row 65
$ws.Range("H65").Value = 4091.7273
$ws.Range("J65").Value = 3882.2
$ws.Range("L65").Value = 19411
$ws.Range("N65").Value = -25651

# Row 107
$ws.Range("H107").Value = 560.32355
$ws.Range("I107").Value = 533.28
$ws.Range("J107").Value = 635.44446
$ws.Range("K107").Value = 533.28
$ws.Range("L107").Value = 635.44446
$ws.Range("M107").Value = 1386.72
$ws.Range("N107").Value = -4475.44446

# Row 109
$ws.Range("H109").Value = 78800
$ws.Range("J109").Value = 78800
$ws.Range("L109").Value = 78800
$ws.Range("N109").Value = -81574

# Row 128
$ws.Range("H128").Value = 75839.5
$ws.Range("J128").Value = 75839.5
$ws.Range("L128").Value = 75839.5
$ws.Range("N128").Value = -85799.5

# Row 138
$ws.Range("H138").Value = 1404374.5
$ws.Range("I138").Value = 3614.2222
$ws.Range("J138").Value = 1629496.8
$ws.Range("K138").Value = 10842.6666
$ws.Range("L138").Value = 4888490.4
$ws.Range("M138").Value = -5702.6666
$ws.Range("N138").Value = -4898770.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 22145.805
$ws.Range("I32").Value = 15819.223
$ws.Range("J32").Value = 25596.666
$ws.Range("K32").Value = 15819.223
$ws.Range("L32").Value = 25596.666
$ws.Range("M32").Value = -15532.223
$ws.Range("N32").Value = -26170.666

# Row 112
$ws.Range("H112").Value = 15104
$ws.Range("J112").Value = 15104
$ws.Range("L112").Value = 15104
$ws.Range("N112").Value = -18058

# Row 115
$ws.Range("H115").Value = 20996.666
$ws.Range("J115").Value = 20996.666
$ws.Range("L115").Value = 20996.666
$ws.Range("N115").Value = -24130.666

$ws = $wb.Worksheets.Item("BSM")
# Row 6
$ws.Range("H6").Value = 40356
$ws.Range("J6").Value = 40356
$ws.Range("L6").Value = 40356
$ws.Range("N6").Value = -40582

# Row 50
$ws.Range("H50").Value = 39800
$ws.Range("J50").Value = 39800
$ws.Range("L50").Value = 39800
$ws.Range("N50").Value = -40948

# Row 53
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

# Row 60
$ws.Range("H60").Value = 49997.145
$ws.Range("J60").Value = 49997.145
$ws.Range("L60").Value = 49997.145
$ws.Range("N60").Value = -51195.145

# Row 104
$ws.Range("H104").Value = 60684
$ws.Range("J104").Value = 60684
$ws.Range("L104").Value = 60684
$ws.Range("N104").Value = -67672

# Row 117
$ws.Range("H117").Value = 79800
$ws.Range("J117").Value = 79800
$ws.Range("L117").Value = 79800
$ws.Range("N117").Value = -88978

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4859.341
$ws.Range("I31").Value = 2136.3333
$ws.Range("J31").Value = 5880.4688
$ws.Range("K31").Value = 2136.3333
$ws.Range("L31").Value = 5880.4688
$ws.Range("M31").Value = -1841.3333
$ws.Range("N31").Value = -6470.4688

# Row 34
$ws.Range("H34").Value = 4859.341
$ws.Range("I34").Value = 2136.3333
$ws.Range("J34").Value = 5880.4688
$ws.Range("K34").Value = 2136.3333
$ws.Range("L34").Value = 5880.4688
$ws.Range("M34").Value = -1934.3333
$ws.Range("N34").Value = -6284.4688

# Row 86
$ws.Range("H86").Value = 2933.3333
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 4800
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 4800
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -7046

# Row 89
$ws.Range("H89").Value = 2933.3333
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 4800
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 24000
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -35232

# Row 118
$ws.Range("H118").Value = 52266.668
$ws.Range("J118").Value = 55900
$ws.Range("L118").Value = 55900
$ws.Range("N118").Value = -59214

# Row 129
$ws.Range("H129").Value = 79800
$ws.Range("J129").Value = 79800
$ws.Range("L129").Value = 79800
$ws.Range("N129").Value = -89800

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 44.842106
$ws.Range("I2").Value = 10.666667
$ws.Range("J2").Value = 55.448277
$ws.Range("K2").Value = 64.00000199999999
$ws.Range("L2").Value = 332.689662
$ws.Range("M2").Value = 48.99999800000001
$ws.Range("N2").Value = -558.689662

# Row 7
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -188
$ws.Range("N7").ClearContents()

# Row 34
$ws.Range("H34").Value = 6463.5713
$ws.Range("I34").Value = 111.25
$ws.Range("J34").Value = 14933.333
$ws.Range("K34").Value = 333.75
$ws.Range("L34").Value = 44799.999
$ws.Range("M34").Value = -249.75
$ws.Range("N34").Value = -44967.999

# Row 39
$ws.Range("H39").Value = 7432.0713
$ws.Range("J39").Value = 7432.0713
$ws.Range("L39").Value = 22296.2139
$ws.Range("N39").Value = -22884.2139

# Row 55
$ws.Range("H55").Value = 5825
$ws.Range("J55").Value = 7666.6665
$ws.Range("L55").Value = 22999.9995
$ws.Range("N55").Value = -23353.9995

# Row 129
$ws.Range("H129").Value = 3720
$ws.Range("J129").Value = 1500
$ws.Range("L129").Value = 4500
$ws.Range("N129").Value = -14500

# Row 131
$ws.Range("H131").Value = 536.9394
$ws.Range("I131").Value = 290.29413
$ws.Range("J131").Value = 799
$ws.Range("K131").Value = 870.88239
$ws.Range("L131").Value = 2397
$ws.Range("M131").Value = 4169.11761
$ws.Range("N131").Value = -12477

$ws = $wb.Worksheets.Item("GSM")
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 4710.9775
$ws.Range("I136").Value = 2953.8708
$ws.Range("J136").Value = 8601.714
$ws.Range("K136").Value = 8861.6124
$ws.Range("L136").Value = 25805.142
$ws.Range("M136").Value = -6311.6124
$ws.Range("N136").Value = -30905.142

# Row 137
$ws.Range("H137").Value = 52311.285
$ws.Range("J137").Value = 54363.168
$ws.Range("L137").Value = 54363.168
$ws.Range("N137").Value = -64563.168

$ws = $wb.Worksheets.Item("WVR")
# Row 121
$ws.Range("H121").Value = 39110
$ws.Range("J121").Value = 39110
$ws.Range("L121").Value = 39110
$ws.Range("N121").Value = -42604

# Row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# Row 132
$ws.Range("H132").Value = 1767.9546
$ws.Range("I132").Value = 1142.24
$ws.Range("J132").Value = 2591.2632
$ws.Range("K132").Value = 3426.72
$ws.Range("L132").Value = 7773.7896
$ws.Range("M132").Value = -896.7200000000003
$ws.Range("N132").Value = -12833.7896
